$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = -0.071696628145104316
$ws.Range("B1").Value = 0.071696627025601201

$ws.Range("A2").Value = 0.057219944403532008
$ws.Range("B2").Value = -0.057219945564075191

$ws.Range("A3").Value = -0.0136755813339177
$ws.Range("B3").Value = 0.01367558009700863
